$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-20 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-21 Friday", 2) | Out-Null
$d.Content.Find.Execute("18×45=810", $true, $false, $false, $false, $false, $true, 1, $false, "49×84=4116", 2) | Out-Null
$d.Content.Find.Execute("32×13=416", $true, $false, $false, $false, $false, $true, 1, $false, "29×62=1798", 2) | Out-Null
$d.Content.Find.Execute("15×66=990", $true, $false, $false, $false, $false, $true, 1, $false, "50×26=1300", 2) | Out-Null
$d.Content.Find.Execute("83×19=1577", $true, $false, $false, $false, $false, $true, 1, $false, "99×11=1089", 2) | Out-Null
$d.Content.Find.Execute("81×93=7533", $true, $false, $false, $false, $false, $true, 1, $false, "60×93=5580", 2) | Out-Null
$d.Content.Find.Execute("25×43=1075", $true, $false, $false, $false, $false, $true, 1, $false, "24×44=1056", 2) | Out-Null
$d.Content.Find.Execute("31×68=2108", $true, $false, $false, $false, $false, $true, 1, $false, "29×97=2813", 2) | Out-Null
$d.Content.Find.Execute("42×13=546", $true, $false, $false, $false, $false, $true, 1, $false, "24×74=1776", 2) | Out-Null
$d.Content.Find.Execute("39×28=1092", $true, $false, $false, $false, $false, $true, 1, $false, "43×38=1634", 2) | Out-Null
$d.Content.Find.Execute("98×14=1372", $true, $false, $false, $false, $false, $true, 1, $false, "29×20=580", 2) | Out-Null
$d.Content.Find.Execute("21×86=1806", $true, $false, $false, $false, $false, $true, 1, $false, "15×55=825", 2) | Out-Null
$d.Content.Find.Execute("78×84=6552", $true, $false, $false, $false, $false, $true, 1, $false, "93×57=5301", 2) | Out-Null
$d.Content.Find.Execute("69×82=5658", $true, $false, $false, $false, $false, $true, 1, $false, "97×18=1746", 2) | Out-Null
$d.Content.Find.Execute("85×67=5695", $true, $false, $false, $false, $false, $true, 1, $false, "14×94=1316", 2) | Out-Null
$d.Content.Find.Execute("61×60=3660", $true, $false, $false, $false, $false, $true, 1, $false, "68×46=3128", 2) | Out-Null
$d.Content.Find.Execute("98×93=9114", $true, $false, $false, $false, $false, $true, 1, $false, "54×96=5184", 2) | Out-Null
$d.Content.Find.Execute("78×17=1326", $true, $false, $false, $false, $false, $true, 1, $false, "70×48=3360", 2) | Out-Null
$d.Content.Find.Execute("44×14=616", $true, $false, $false, $false, $false, $true, 1, $false, "16×86=1376", 2) | Out-Null
$d.Content.Find.Execute("62×95=5890", $true, $false, $false, $false, $false, $true, 1, $false, "52×87=4524", 2) | Out-Null
$d.Content.Find.Execute("24×97=2328", $true, $false, $false, $false, $false, $true, 1, $false, "91×60=5460", 2) | Out-Null
$d.Content.Find.Execute("64×18=1152", $true, $false, $false, $false, $false, $true, 1, $false, "30×33=990", 2) | Out-Null
$d.Content.Find.Execute("47×27=1269", $true, $false, $false, $false, $false, $true, 1, $false, "22×80=1760", 2) | Out-Null
$d.Content.Find.Execute("60×81=4860", $true, $false, $false, $false, $false, $true, 1, $false, "43×43=1849", 2) | Out-Null
$d.Content.Find.Execute("23×73=1679", $true, $false, $false, $false, $false, $true, 1, $false, "81×61=4941", 2) | Out-Null
$d.Content.Find.Execute("87×73=6351", $true, $false, $false, $false, $false, $true, 1, $false, "98×48=4704", 2) | Out-Null
